$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Single-cell value corrections in column C (rows 860-992) ---
$ws.Cells.Item(860, 3).Value = 20
$ws.Cells.Item(870, 3).Value = 16
$ws.Cells.Item(913, 3).Value = 23
$ws.Cells.Item(917, 3).Value = 17
$ws.Cells.Item(924, 3).Value = 5
$ws.Cells.Item(946, 3).Value = 24
$ws.Cells.Item(951, 3).Value = 30
$ws.Cells.Item(956, 3).Value = 34
$ws.Cells.Item(965, 3).Value = 11
$ws.Cells.Item(967, 3).Value = 29
$ws.Cells.Item(968, 3).Value = 2
$ws.Cells.Item(970, 3).Value = 13
$ws.Cells.Item(971, 3).Value = 29
$ws.Cells.Item(975, 3).Value = 31
$ws.Cells.Item(979, 3).Value = 15
$ws.Cells.Item(980, 3).Value = 33
$ws.Cells.Item(983, 3).Value = 34
$ws.Cells.Item(985, 3).Value = 9
$ws.Cells.Item(987, 3).Value = 31
$ws.Cells.Item(991, 3).Value = 12
$ws.Cells.Item(992, 3).Value = 26

# --- Rewrite rows 993-1016: date/age-group labels shift down, new rows appended ---
$ws.Cells.Item(993, 1).Value = 44161
$ws.Cells.Item(993, 2).Value = "40-49"
$ws.Cells.Item(993, 3).Value = 2
$ws.Cells.Item(994, 1).Value = 44161
$ws.Cells.Item(994, 2).Value = "50-59"
$ws.Cells.Item(994, 3).Value = 4
$ws.Cells.Item(995, 1).Value = 44161
$ws.Cells.Item(995, 2).Value = "60-69"
$ws.Cells.Item(995, 3).Value = 9
$ws.Cells.Item(996, 1).Value = 44161
$ws.Cells.Item(996, 2).Value = "70-79"
$ws.Cells.Item(996, 3).Value = 13
$ws.Cells.Item(997, 1).Value = 44161
$ws.Cells.Item(997, 2).Value = "80+"
$ws.Cells.Item(997, 3).Value = 33
$ws.Cells.Item(998, 1).Value = 44162
$ws.Cells.Item(998, 2).Value = "50-59"
$ws.Cells.Item(998, 3).Value = 2
$ws.Cells.Item(999, 1).Value = 44162
$ws.Cells.Item(999, 2).Value = "60-69"
$ws.Cells.Item(999, 3).Value = 7
$ws.Cells.Item(1000, 1).Value = 44162
$ws.Cells.Item(1000, 2).Value = "70-79"
$ws.Cells.Item(1000, 3).Value = 26
$ws.Cells.Item(1001, 1).Value = 44162
$ws.Cells.Item(1001, 2).Value = "80+"
$ws.Cells.Item(1001, 3).Value = 25
$ws.Cells.Item(1002, 1).Value = 44163
$ws.Cells.Item(1002, 2).Value = "30-39"
$ws.Cells.Item(1002, 3).Value = 2
$ws.Cells.Item(1003, 1).Value = 44163
$ws.Cells.Item(1003, 2).Value = "40-49"
$ws.Cells.Item(1003, 3).Value = 1
$ws.Cells.Item(1004, 1).Value = 44163
$ws.Cells.Item(1004, 2).Value = "50-59"
$ws.Cells.Item(1004, 3).Value = 2
$ws.Cells.Item(1005, 1).Value = 44163
$ws.Cells.Item(1005, 2).Value = "60-69"
$ws.Cells.Item(1005, 3).Value = 8
$ws.Cells.Item(1006, 1).Value = 44163
$ws.Cells.Item(1006, 2).Value = "70-79"
$ws.Cells.Item(1006, 3).Value = 9
$ws.Cells.Item(1007, 1).Value = 44163
$ws.Cells.Item(1007, 2).Value = "80+"
$ws.Cells.Item(1007, 3).Value = 19
$ws.Cells.Item(1008, 1).Value = 44164
$ws.Cells.Item(1008, 2).Value = "50-59"
$ws.Cells.Item(1008, 3).Value = 1
$ws.Cells.Item(1009, 1).Value = 44164
$ws.Cells.Item(1009, 2).Value = "60-69"
$ws.Cells.Item(1009, 3).Value = 3
$ws.Cells.Item(1010, 1).Value = 44164
$ws.Cells.Item(1010, 2).Value = "70-79"
$ws.Cells.Item(1010, 3).Value = 15
$ws.Cells.Item(1011, 1).Value = 44164
$ws.Cells.Item(1011, 2).Value = "80+"
$ws.Cells.Item(1011, 3).Value = 34
$ws.Cells.Item(1012, 1).Value = 44165
$ws.Cells.Item(1012, 2).Value = "0-19"
$ws.Cells.Item(1012, 3).Value = 1
$ws.Cells.Item(1013, 1).Value = 44165
$ws.Cells.Item(1013, 2).Value = "50-59"
$ws.Cells.Item(1013, 3).Value = 2
$ws.Cells.Item(1014, 1).Value = 44165
$ws.Cells.Item(1014, 2).Value = "60-69"
$ws.Cells.Item(1014, 3).Value = 2
$ws.Cells.Item(1015, 1).Value = 44165
$ws.Cells.Item(1015, 2).Value = "70-79"
$ws.Cells.Item(1015, 3).Value = 5
$ws.Cells.Item(1016, 1).Value = 44165
$ws.Cells.Item(1016, 2).Value = "80+"
$ws.Cells.Item(1016, 3).Value = 15

# --- Apply the date number format to column A for the newly added rows ---
$dateFormat = $ws.Cells.Item(992, 1).NumberFormat
$ws.Range("A1011:A1016").NumberFormat = $dateFormat
